$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

$ws_ALC.Range("H17").Value = 1638.7742
$ws_ALC.Range("J17").Value = 1324.6333
$ws_ALC.Range("L17").Value = 3973.8999
$ws_ALC.Range("N17").Value = -4309.8999
$ws_ALC.Range("I40").Value = 2495
$ws_ALC.Range("J40").Value = 3401
$ws_ALC.Range("K40").Value = 2495
$ws_ALC.Range("L40").Value = 3401
$ws_ALC.Range("M40").Value = -2320
$ws_ALC.Range("N40").Value = -3751
$ws_ALC.Range("H80").Value = 1291.0416
$ws_ALC.Range("J80").Value = 647.75
$ws_ALC.Range("L80").Value = 1943.25
$ws_ALC.Range("N80").Value = -3939.25
$ws_ALC.Range("H83").Value = 1291.0416
$ws_ALC.Range("J83").Value = 647.75
$ws_ALC.Range("L83").Value = 5829.75
$ws_ALC.Range("N83").Value = -15813.75
$ws_ALC.Range("H93").Value = 84600
$ws_ALC.Range("J93").Value = 84600
$ws_ALC.Range("L93").Value = 84600
$ws_ALC.Range("N93").Value = -89592
$ws_ALC.Range("H112").Value = 7857.143
$ws_ALC.Range("J112").Value = 8983.333000000001
$ws_ALC.Range("L112").Value = 26949.999
$ws_ALC.Range("N112").Value = -29165.999
$ws_ALC.Range("H125").Value = 1118.1428
$ws_ALC.Range("J125").Value = 1037.5
$ws_ALC.Range("L125").Value = 9337.5
$ws_ALC.Range("N125").Value = -14257.5
$ws_ALC.Range("H138").Value = 1734.836
$ws_ALC.Range("J138").Value = 2185.6875
$ws_ALC.Range("L138").Value = 6557.0625
$ws_ALC.Range("N138").Value = -16837.0625
$ws_ARM.Range("H32").Value = 3775.1355
$ws_ARM.Range("I32").Value = 2503.4583
$ws_ARM.Range("J32").Value = 9324.272000000001
$ws_ARM.Range("K32").Value = 2503.4583
$ws_ARM.Range("L32").Value = 9324.272000000001
$ws_ARM.Range("M32").Value = -2216.4583
$ws_ARM.Range("N32").Value = -9898.272000000001
$ws_ARM.Range("H53").Value = 20000
$ws_ARM.Range("I53").Value = 0
$ws_ARM.Range("K53").Value = 0
$ws_ARM.Range("M53").ClearContents()
$ws_ARM.Range("H61").Value = 3556.818
$ws_ARM.Range("I61").Value = 2614.7058
$ws_ARM.Range("K61").Value = 2614.7058
$ws_ARM.Range("M61").Value = -2402.7058
$ws_ARM.Range("H110").Value = 1518.0416
$ws_ARM.Range("I110").Value = 1148.6957
$ws_ARM.Range("K110").Value = 1148.6957
$ws_ARM.Range("M110").Value = 896.3043
$ws_ARM.Range("H136").Value = 3556.818
$ws_ARM.Range("I136").Value = 2614.7058
$ws_ARM.Range("K136").Value = 7844.117400000001
$ws_ARM.Range("M136").Value = -5294.117400000001
$ws_BSM.Range("H62").Value = 10000
$ws_BSM.Range("I62").Value = 10000
$ws_BSM.Range("K62").Value = 10000
$ws_BSM.Range("M62").Value = -9314
$ws_BSM.Range("H65").Value = 10000
$ws_BSM.Range("I65").Value = 10000
$ws_BSM.Range("K65").Value = 30000
$ws_BSM.Range("M65").Value = -26568
$ws_BSM.Range("H134").Value = 4020.1052
$ws_BSM.Range("I134").Value = 3471.2222
$ws_BSM.Range("J134").Value = 13900
$ws_BSM.Range("K134").Value = 10413.6666
$ws_BSM.Range("L134").Value = 41700
$ws_BSM.Range("M134").Value = -7878.6666
$ws_BSM.Range("N134").Value = -46770
$ws_CRP.Range("H63").Value = 30000
$ws_CRP.Range("J63").Value = 30000
$ws_CRP.Range("L63").Value = 30000
$ws_CRP.Range("N63").Value = -31372
$ws_CRP.Range("H66").Value = 30000
$ws_CRP.Range("J66").Value = 30000
$ws_CRP.Range("L66").Value = 90000
$ws_CRP.Range("N66").Value = -96864
$ws_CRP.Range("H99").Value = 2849.875
$ws_CRP.Range("I99").Value = 2500
$ws_CRP.Range("K99").Value = 2500
$ws_CRP.Range("M99").Value = -1002
$ws_CRP.Range("H126").Value = 2849.875
$ws_CRP.Range("I126").Value = 2500
$ws_CRP.Range("K126").Value = 7500
$ws_CRP.Range("M126").Value = -5030
$ws_CRP.Range("H134").Value = 1186.0646
$ws_CRP.Range("I134").Value = 1186.0646
$ws_CRP.Range("J134").Value = 0
$ws_CRP.Range("K134").Value = 3558.1938
$ws_CRP.Range("L134").Value = 0
$ws_CRP.Range("M134").Value = -1023.1938
$ws_CRP.Range("N134").ClearContents()
$ws_CUL.Range("H108").Value = 3000.5
$ws_CUL.Range("I108").Value = 3000.5
$ws_CUL.Range("K108").Value = 9001.5
$ws_CUL.Range("M108").Value = -6121.5
$ws_CUL.Range("H113").Value = 944.6667
$ws_CUL.Range("I113").Value = 1000
$ws_CUL.Range("J113").Value = 939.63635
$ws_CUL.Range("K113").Value = 3000
$ws_CUL.Range("L113").Value = 2818.90905
$ws_CUL.Range("M113").Value = -830
$ws_CUL.Range("N113").Value = -7158.90905
$ws_CUL.Range("H121").Value = 820.375
$ws_CUL.Range("J121").Value = 822.1667
$ws_CUL.Range("L121").Value = 2466.5001
$ws_CUL.Range("N121").Value = -5086.5001
$ws_CUL.Range("H131").Value = 784.39795
$ws_CUL.Range("J131").Value = 813.30334
$ws_CUL.Range("L131").Value = 2439.91002
$ws_CUL.Range("N131").Value = -12519.91002
$ws_CUL.Range("H134").Value = 2580.4666
$ws_CUL.Range("I134").Value = 2065
$ws_CUL.Range("J134").Value = 3998
$ws_CUL.Range("K134").Value = 6195
$ws_CUL.Range("L134").Value = 11994
$ws_CUL.Range("M134").Value = -1125
$ws_CUL.Range("N134").Value = -22134
$ws_GSM.Range("H113").Value = 1021.9091
$ws_GSM.Range("I113").Value = 719.6
$ws_GSM.Range("J113").Value = 1273.8334
$ws_GSM.Range("K113").Value = 719.6
$ws_GSM.Range("L113").Value = 1273.8334
$ws_GSM.Range("M113").Value = 1450.4
$ws_GSM.Range("N113").Value = -5613.8334
$ws_LTW.Range("H7").Value = 4279.0835
$ws_LTW.Range("I7").Value = 2373.889
$ws_LTW.Range("J7").Value = 9994.666999999999
$ws_LTW.Range("K7").Value = 2373.889
$ws_LTW.Range("L7").Value = 9994.666999999999
$ws_LTW.Range("M7").Value = -2261.889
$ws_LTW.Range("N7").Value = -10218.667
$ws_LTW.Range("H46").Value = 3574.8
$ws_LTW.Range("I46").Value = 0
$ws_LTW.Range("J46").Value = 3574.8
$ws_LTW.Range("K46").Value = 0
$ws_LTW.Range("L46").Value = 3574.8
$ws_LTW.Range("M46").ClearContents()
$ws_LTW.Range("N46").Value = -3950.8
$ws_LTW.Range("H61").Value = 3699.375
$ws_LTW.Range("J61").Value = 4999.25
$ws_LTW.Range("L61").Value = 4999.25
$ws_LTW.Range("N61").Value = -5403.25
$ws_LTW.Range("H113").Value = 3699.375
$ws_LTW.Range("J113").Value = 4999.25
$ws_LTW.Range("L113").Value = 4999.25
$ws_LTW.Range("N113").Value = -9339.25
$ws_LTW.Range("H126").Value = 4279.0835
$ws_LTW.Range("I126").Value = 2373.889
$ws_LTW.Range("J126").Value = 9994.666999999999
$ws_LTW.Range("K126").Value = 7121.667
$ws_LTW.Range("L126").Value = 29984.001
$ws_LTW.Range("M126").Value = -4651.667
$ws_LTW.Range("N126").Value = -34924.001
$ws_LTW.Range("H132").Value = 3402.5
$ws_LTW.Range("I132").Value = 2654.4167
$ws_LTW.Range("K132").Value = 7963.250100000001
$ws_LTW.Range("M132").Value = -5433.250100000001
$ws_LTW.Range("H136").Value = 3773.261
$ws_LTW.Range("I136").Value = 2785.6667
$ws_LTW.Range("K136").Value = 8357.000100000001
$ws_LTW.Range("M136").Value = -5807.000100000001
$ws_WVR.Range("H81").Value = 2299.8
$ws_WVR.Range("I81").Value = 1999.5
$ws_WVR.Range("K81").Value = 3999
$ws_WVR.Range("M81").Value = -2938
$ws_WVR.Range("H84").Value = 2299.8
$ws_WVR.Range("I84").Value = 1999.5
$ws_WVR.Range("K84").Value = 19995
$ws_WVR.Range("M84").Value = -14691
$ws_WVR.Range("H132").Value = 1727.0344
$ws_WVR.Range("I132").Value = 1180.35
$ws_WVR.Range("J132").Value = 2941.889
$ws_WVR.Range("K132").Value = 3541.05
$ws_WVR.Range("L132").Value = 8825.667000000001
$ws_WVR.Range("M132").Value = -1011.05
$ws_WVR.Range("N132").Value = -13885.667

Write-Host "Applied all updates"